$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (and any touched E cells) to Text format before assignment
# so numeric-looking strings like "1.00" / "544.25" are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '60.755.87'
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").Value = '2.376.46'
$ws.Range("E3").Value = '  -3.04%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '544.25'
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("D6").Value = '141.15'
$ws.Range("E6").Value = '  -2.12%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '0.535'
$ws.Range("E8").Value = '  -9.95%  '
$ws.Range("D9").Value = '2.375.60'
$ws.Range("E9").Value = '  -2.96%  '
$ws.Range("E10").Value = '  -1.09%  '
$ws.Range("E11").Value = '  +0.70%  '
$ws.Range("D12").Value = '5.34'
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("D13").Value = '0.344'
$ws.Range("E13").Value = '  -1.62%  '
$ws.Range("D14").Value = '25.47'
$ws.Range("E14").Value = '  -1.54%  '
$ws.Range("D15").Value = '2.802.96'
$ws.Range("E15").Value = '  -3.05%  '
$ws.Range("D16").Value = '0.0000163'
$ws.Range("E16").Value = '  +0.45%  '
$ws.Range("D17").Value = '60.636.05'
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("D18").Value = '2.374.83'
$ws.Range("E18").Value = '  -3.42%  '
$ws.Range("D19").Value = '10.61'
$ws.Range("E19").Value = '  -3.51%  '
$ws.Range("D20").Value = '4.09'
$ws.Range("E20").Value = '  -1.42%  '
$ws.Range("D21").Value = '315.76'
$ws.Range("E21").Value = '  -0.53%  '
$ws.Range("D22").Value = '6.70'
$ws.Range("E22").Value = '  -2.71%  '
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("E24").Value = '  +3.17%  '
$ws.Range("D25").Value = '62.72'
$ws.Range("E25").Value = '  -0.83%  '
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("D27").Value = '2.489.16'
$ws.Range("E27").Value = '  -3.47%  '
$ws.Range("D28").Value = '0.0₃0927'
$ws.Range("E28").Value = '  -4.10%  '
$ws.Range("D29").Value = '7.72'
$ws.Range("E29").Value = '  +2.66%  '
$ws.Range("D30").Value = '521.75'
$ws.Range("E30").Value = '  -2.02%  '
$ws.Range("E31").Value = '  -3.28%  '
$ws.Range("D32").Value = '7.97'
$ws.Range("E32").Value = '  -3.83%  '
$ws.Range("E33").Value = '  -2.82%  '
$ws.Range("D34").Value = '1.83'
$ws.Range("E34").Value = '  -2.78%  '
$ws.Range("E35").Value = '  -0.32%  '
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("D37").Value = '5.47'
$ws.Range("E37").Value = '  -5.84%  '
$ws.Range("D38").Value = '4.65'
$ws.Range("E38").Value = '  -2.93%  '
$ws.Range("D39").Value = '0.375'
$ws.Range("E39").Value = '  +0.09%  '
$ws.Range("D40").Value = '17.96'
$ws.Range("E40").Value = '  -2.31%  '
$ws.Range("D41").Value = '1.72'
$ws.Range("E41").Value = '  +1.91%  '
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("D43").Value = '136.75'
$ws.Range("E43").Value = '  -6.09%  '
$ws.Range("D44").Value = '40.26'
$ws.Range("E44").Value = '  +1.22%  '
$ws.Range("D45").Value = '2.23'
$ws.Range("E45").Value = '  -0.95%  '
$ws.Range("D46").Value = '139.50'
$ws.Range("E46").Value = '  -4.41%  '
$ws.Range("D47").Value = '3.54'
$ws.Range("E47").Value = '  +0.21%  '
$ws.Range("D48").Value = '20.29'
$ws.Range("E48").Value = '  -1.97%  '
$ws.Range("D49").Value = '0.0516'
$ws.Range("E49").Value = '  -2.00%  '
$ws.Range("D50").Value = '0.575'
$ws.Range("E50").Value = '  -1.00%  '
$ws.Range("D51").Value = '0.0911'
$ws.Range("E51").Value = '  -2.60%  '

# Restore original (default/General) formatting so no stray style indices are introduced,
# matching the source workbook where these cells carry no explicit style.
$ws.Range("D2:D51").ClearFormats()
